$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need NumberFormat forced to
# Text first, otherwise Excel auto-converts the literal price string (losing
# trailing zeros / thousands-dot formatting) the same way it would if you typed
# "1.00" into a General-formatted cell.
$textCells = @("D5", "D6", "D11", "D12", "D13", "D14", "D19", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D32", "D33", "D34", "D40", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = '69.332.43'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").Value = '3.680.97'
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '682.73'
$ws.Range("E5").Value = '  -3.77%  '
$ws.Range("D6").Value = '162.68'
$ws.Range("E6").Value = '  -4.62%  '
$ws.Range("D7").Value = '3.680.71'
$ws.Range("E7").Value = '  -3.37%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E10").Value = '  -7.60%  '
$ws.Range("D11").Value = '7.26'
$ws.Range("E11").Value = '  -3.59%  '
$ws.Range("D12").Value = '0.450'
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("D13").Value = '0.0000237'
$ws.Range("E13").Value = '  -6.10%  '
$ws.Range("D14").Value = '33.61'
$ws.Range("E14").Value = '  -6.61%  '
$ws.Range("D15").Value = '4.304.11'
$ws.Range("E15").Value = '  -3.34%  '
$ws.Range("D16").Value = '3.683.76'
$ws.Range("E16").Value = '  -3.41%  '
$ws.Range("D17").Value = '69.340.95'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("E18").Value = '  -1.80%  '
$ws.Range("D19").Value = '16.33'
$ws.Range("E19").Value = '  -5.66%  '
$ws.Range("E20").Value = '  -6.89%  '
$ws.Range("D21").Value = '479.40'
$ws.Range("E21").Value = '  -3.35%  '
$ws.Range("D22").Value = '9.81'
$ws.Range("E22").Value = '  -7.86%  '
$ws.Range("D23").Value = '0.667'
$ws.Range("D24").Value = '80.13'
$ws.Range("D25").Value = '3.830.10'
$ws.Range("E25").Value = '  -3.32%  '
$ws.Range("E26").Value = '  -11.01%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '11.52'
$ws.Range("E27").Value = '  -4.71%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '9.60'
$ws.Range("E29").Value = '  -7.71%  '
$ws.Range("E30").Value = '  -10.50%  '
$ws.Range("E31").Value = '  -10.98%  '
$ws.Range("D32").Value = '2.11'
$ws.Range("E32").Value = '  -5.70%  '
$ws.Range("D33").Value = '6.87'
$ws.Range("E33").Value = '  -6.36%  '
$ws.Range("D34").Value = '27.10'
$ws.Range("E34").Value = '  -6.87%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("E36").Value = '  -4.59%  '
$ws.Range("D37").Value = '3.644.97'
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("E38").Value = '  -5.72%  '
$ws.Range("E39").Value = '  +3.45%  '
$ws.Range("D40").Value = '0.0943'
$ws.Range("E40").Value = '  -7.14%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  -6.15%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '0.959'
$ws.Range("E44").Value = '  -7.62%  '
$ws.Range("D45").Value = '48.17'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").Value = '154.47'
$ws.Range("E46").Value = '  -7.00%  '
$ws.Range("D47").Value = '2.85'
$ws.Range("E47").Value = '  -12.04%  '
$ws.Range("D48").Value = '1.34'
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").Value = '0.000281'
$ws.Range("E49").Value = '  -12.22%  '
$ws.Range("D50").Value = '391.97'
$ws.Range("E50").Value = '  -7.79%  '
$ws.Range("E51").Value = '  -5.86%  '
